$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two blank rows after row 4 (they become rows 5 and 6) so the
#    "Warm-up questions" block gets two extra sub-question rows.
# ---------------------------------------------------------------------------
$ws.Range("A5:A6").EntireRow.Insert()

# ---------------------------------------------------------------------------
# 2. Set row heights for the whole table (values taken from the target file).
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30.55
$ws.Rows.Item(5).RowHeight = 17.1
$ws.Rows.Item(6).RowHeight = 19.45
$ws.Rows.Item(7).RowHeight = 24.05
$ws.Rows.Item(8).RowHeight = 35.5
$ws.Rows.Item(9).RowHeight = 31.3
$ws.Rows.Item(10).RowHeight = 78.3
$ws.Rows.Item(11).RowHeight = 17.35
$ws.Rows.Item(12).RowHeight = 17.35

# ---------------------------------------------------------------------------
# 3. Copy formatting onto the newly inserted / newly used cells so that the
#    look (font, borders, wrap) matches the rest of the table.
#    - Style "A4-like" (bordered, wrapping, Times New Roman 14) -> A/B content cells
#    - Style "C3-like" (no border, wrapping, Times New Roman 10) -> C column notes
# ---------------------------------------------------------------------------
$ws.Range("A4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A5:B6").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C3").Copy()
$ws.Range("C4:C12").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# The engine's OOXML writer does not always round-trip a pre-existing
# "wrapText" flag faithfully, so make sure every cell that must wrap its
# text has that flag set explicitly (this does not touch the cells that
# should stay non-wrapping: B7:B9, A11:B11).
$ws.Range("A1:B6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A8").WrapText = $true
$ws.Range("A9").WrapText = $true
$ws.Range("A10:B10").WrapText = $true
$ws.Range("C2:C12").WrapText = $true

# ---------------------------------------------------------------------------
# 4. Fill in the new / changed cell values.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Would you introduce yourself? "
$ws.Range("C4").Value = "Establish rapport with participants by asking them some general information questions like " + [char]34 + "Would you introduce yourself?" + [char]34 + " and " + [char]34 + "What is your occupation?" + [char]34 + " " + [char]34 + "What type(s) of mobile device(s) do you use?" + [char]34 + " etc. "

$ws.Range("B5").Value = "What is your occupation? "

$ws.Range("B6").Value = "What type of mobile device do you use?"

$ws.Range("C7").Value = "Give the participants an overview of what they are looking at and ask general questions. For example, you could ask: What do you think the screen is for? "

$ws.Range("C8").Value = "Let participants know how you will guide them through the tasks. Give an overview of the tasks they are required to do. Remember, the test is on the navigation and functions offered by the screen. Write down the answers. "

$ws.Range("C9").Value = "After each task, ask a few questions and conduct a rating exercise to elicit more feedback and information about the users' reasoning. "

$ws.Range("B10").Value = "Finaly, I would like to ask you are there anything that you feel to be missing in our web page? Thank you very mach for your participation it this test!"
$ws.Range("C10").Value = "You don't have to ask many questions, but it's essential to have a few prepared. For example, an additional question could be: Was there anything you felt was missing or thought could have been better? Answer any questions the participants may have. Thank the participants."

# ---------------------------------------------------------------------------
# 5. Re-merge the trailing blank row (it moved from row 9 to row 11).
# ---------------------------------------------------------------------------
$ws.Range("A11:B11").Merge()

# ---------------------------------------------------------------------------
# 6. Update the selection to match the target file (row 5 selected).
# ---------------------------------------------------------------------------
$ws.Range("A5").EntireRow.Select()
